# "small change for tracking dataset"
# Adds new experiment-log rows (10_2 .. 12_4) to Sheet2, clears the now-
# superseded K91 value, and updates the saved selection / active sheet
# view state on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet2 ("exp log") data additions
# ---------------------------------------------------------------------

# Row 91 gains a date / note; row 92's old K value (0.02) is superseded
# because row 91 no longer carries a lambda of its own.
$ws2.Range("B91").Value = "Feb_20"
$ws2.Range("C91").Value = "supervised"
$ws2.Range("F91").Value = 10000
$ws2.Range("K91").ClearContents()

$ws2.Range("A92").Value = "10_2"
$ws2.Range("F92").Value = 50000
$ws2.Range("G92").Value = "all"
$ws2.Range("K92").Value = 0.02
$ws2.Range("O92").Value = 0.001
$ws2.Range("P92").Value = "0.07, 0.56, 1.21"

$ws2.Range("A93").Value = "10_3"
$ws2.Range("K93").Value = 0.05
$ws2.Range("O93").Value = 0.001
$ws2.Range("P93").Value = "0.12, 0.61, 0.85"

$ws2.Range("A94").Value = "10_4"
$ws2.Range("K94").Value = 0.1
$ws2.Range("O94").Value = 0.0005
$ws2.Range("P94").Value = "0.19, 0.46, 0.64"

$ws2.Range("A95").Value = "10_5"
$ws2.Range("C95").Value = "Supervised-pretrain"
$ws2.Range("K95").Value = 0.02
$ws2.Range("O95").Value = 0.001
$ws2.Range("P95").Value = "0.04, 0.49, 0.55"
$ws2.Range("Q95").Value = "supervised pretrain help a little"

$ws2.Range("A97").Value = "10_4_2"
$ws2.Range("B97").Value = "Feb_21"
$ws2.Range("C97").Value = "Continue 10_4"
$ws2.Range("F97").Value = 20000
$ws2.Range("K97").Value = 0.1
$ws2.Range("O97").Value = 0.0005

$ws2.Range("A98").Value = "10_5_2"
$ws2.Range("C98").Value = "Continue 10_2 with lam=0.05"
$ws2.Range("F98").Value = 50000
$ws2.Range("K98").Value = 0.05
$ws2.Range("O98").Value = 0.0005

$ws2.Range("A100").Value = "10_6"
$ws2.Range("C100").Value = "add more labeled data, from 10_1"
$ws2.Range("F100").Value = 50000
$ws2.Range("K100").Value = 0.03
$ws2.Range("O100").Value = 0.0005
$ws2.Range("P100").Value = "0.12, 0.53, 1.62"

$ws2.Range("A101").Value = "10_7"
$ws2.Range("F101").Value = 50000
$ws2.Range("K101").Value = 0.1
$ws2.Range("O101").Value = 0.0005
$ws2.Range("P101").Value = "0.34, 0.53, 0.71"

$ws2.Range("A103").Value = "10_8"
$ws2.Range("K103").Value = 0.02
$ws2.Range("O103").Value = 0.001
$ws2.Range("P103").Value = "0.10, 0.50, 1.88"
$ws2.Range("Q103").Value = "0.07, 0.48, 1.26"
$ws2.Range("R103").Value = "0.04, 0.50, 0.62"

$ws2.Range("A104").Value = "10_9"
$ws2.Range("K104").Value = 0.05
$ws2.Range("O104").Value = 0.001
$ws2.Range("P104").Value = "0.18, 0.57, 1.45"
$ws2.Range("Q104").Value = "0.11, 0.57, 0.95"

$ws2.Range("A105").Value = "10_10"
$ws2.Range("K105").Value = 0.1
$ws2.Range("O105").Value = 0.001
$ws2.Range("P105").Value = "0.28, 0.66, 0.91"
$ws2.Range("Q105").Value = "0.17, 0.59, 0.60"

$ws2.Range("A106").Value = "10_11"
$ws2.Range("K106").Value = 0.2
$ws2.Range("O106").Value = 0.001
$ws2.Range("P106").Value = "0.33, 0.76, 0.25"
$ws2.Range("Q106").Value = "0.28, 0.72, 0.28"

$ws2.Range("A108").Value = "rescale augmentation"

$ws2.Range("A109").Value = "11_1"
$ws2.Range("G109").Value = "all"
$ws2.Range("K109").Value = 0.02

$ws2.Range("A113").Value = "12_1"
$ws2.Range("G113").Value = "all"
$ws2.Range("K113").Value = 0.002
$ws2.Range("O113").Value = 0.001

$ws2.Range("A114").Value = "12_2"
$ws2.Range("K114").Value = 0.005
$ws2.Range("O114").Value = 0.0005

$ws2.Range("A115").Value = "12_3"
$ws2.Range("K115").Value = 0.001
$ws2.Range("O115").Value = 0.001

$ws2.Range("A116").Value = "12_4"
$ws2.Range("K116").Value = 0.0005
$ws2.Range("O116").Value = 0.001

$ws2.Range("K117").Value = 0.0005
$ws2.Range("O117").Value = 0.0005

$ws2.Range("K118").Value = 0.0001
$ws2.Range("O118").Value = 0.0005

$ws2.Range("K119").Value = 0.0005
$ws2.Range("O119").Value = 0.0001

# ---------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------

# Sheet1: scroll position moves down; saved selection is unchanged (A59).
$ws1.Activate()
$ws1.Range("A59").Select()

# Sheet2: stays the active/tabbed sheet; scroll + selection move to the
# newly-added rows at the bottom of the log.
$ws2.Activate()
$ws2.Range("L120").Select()
